$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logical Operators")

# --- Pass/Fail (column D) ---
$ws.Range("D2").Formula = '=IF(C2>=60,"PASS","FAIL")'
$ws.Range("D3:D16").Formula = '=IF(C3>=60,"PASS","FAIL")'

# --- Outlier (column F) ---
$ws.Range("F2").Formula = '=IF(OR(C2>90,C2<60),"OUTLIER","AVG")'
$ws.Range("F3:F16").Formula = '=IF(OR(C3>90,C3<60),"OUTLIER","AVG")'

# --- Award (column G) ---
$ws.Range("G2").Formula = '=IF(AND(B2="M",C2>95),"Male Achiever",IF(AND(B2="F",C2>95),"Female Achiever","NONE"))'
$ws.Range("G3:G16").Formula = '=IF(AND(B3="M",C3>95),"Male Achiever",IF(AND(B3="F",C3>95),"Female Achiever","NONE"))'

# --- Letter (column E) ---
$ws.Range("E2").Formula = '=IF(C2>=90,"A",IF(AND(C2<=89,C2>=80),"B",IF(AND(C2<=79,C2>=70),"C",IF(AND(C2<=69,C2>=60),"D","F"))))'
$ws.Range("E3").Formula = '=IF(C3>=90,"A",IF(AND(C3<=89,C3>=80),"B",IF(AND(C3<=79,C3>=70),"C",IF(AND(C3<=69,C3>=60),"D","F"))))'
$ws.Range("E4").Formula = '=IF(C4>=90,"A",IF(AND(C4<=89,C4>=80),"B",IF(AND(C4<=79,C4>=70),"C",IF(AND(C4<=69,C4>=60),"D","F"))))'
$ws.Range("E5").Formula = '=IF(C5>=90,"A",IF(AND(C5<=89,C5>=80),"B",IF(AND(C5<=79,C5>=70),"C",IF(AND(C5<=69,C5>=60),"D","F"))))'
$ws.Range("E6").Formula = '=IF(C6>=90,"A",IF(AND(C6<=89,C6>=80),"B",IF(AND(C6<=79,C6>=70),"C",IF(AND(C6<=69,C6>=60),"D","F"))))'
$ws.Range("E7").Formula = '=IF(C7>=90,"A",IF(AND(C7<=89,C7>=80),"B",IF(AND(C7<=79,C7>=70),"C",IF(AND(C7<=69,C7>=60),"D","F"))))'
$ws.Range("E8:E16").Formula = '=IF(C8>=90,"A",IF(AND(C8<=89,C8>=80),"B",IF(AND(C8<=79,C8>=70),"C",IF(AND(C8<=69,C8>=60),"D","F"))))'

# --- View state: make "Logical Operators" the active tab, zoomed to 120%, selection on G22 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("G22").Select()
